$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15 (pushes old rows 15,16 down to 16,17)
$ws.Rows("15:15").Insert()

# Copy the formatting (styles/borders/fills/fonts) from row 14 onto the
# newly inserted row 15 so the new product row matches the others.
$ws.Range("A14:N14").Copy()
$ws.Range("A15:N15").PasteSpecial(-4122)

# Match the row height used by the other product rows
$ws.Rows("15:15").RowHeight = 25.5

# Re-create the merges for the new row 15 (same layout as rows 4-14)
$ws.Range("B15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()

# Fill in the new product row's values
$ws.Range("A15").Value = 12
$ws.Range("B15").Value = "فرد شعر هير كونترول"
$ws.Range("H15").Value = "1:0"
$ws.Range("L15").Value = 180
$ws.Range("N15").Value = "1:0"

# Update the running total row (old row15 -> now row16) to include the new item
$ws.Range("K16").Value = 725

Write-Host "Done"
